$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.337.75'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.407.51'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.64'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.97'
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '2.419.07'
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E10").Value = '  +2.37%  '
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.22'
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.55'
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("D16").Value = '2.847.06'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '61.148.57'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '2.419.33'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.06'
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.71'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.36'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").Value = '  +1.44%  '
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("E24").Value = '  +7.43%  '
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.19'
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '613.93'
$ws.Range("E27").Value = '  +6.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.28'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = '0.0₃0951'
$ws.Range("E29").Value = '  +3.32%  '
$ws.Range("D30").Value = '2.518.99'
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.02'
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").Value = '  +4.02%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.47'
$ws.Range("E35").Value = '  +5.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.67%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.50'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.373'
$ws.Range("E38").Value = '  +1.62%  '
$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.64'
$ws.Range("E39").Value = '  +1.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.32'
$ws.Range("E40").Value = '  +4.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.37'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.56'
$ws.Range("E42").Value = '  +9.70%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.71'
$ws.Range("E43").Value = '  +3.43%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.97'
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("D46").Value = '0.0₆0285'
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.16'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.97'
$ws.Range("E49").Value = '  +3.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.593'
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0511'
$ws.Range("E51").Value = '  +2.38%  '
